$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 38 (shifts the existing data for this
# product's "2021-11-19" block into rows 40.. and opens up 38:39 for the
# newly-reported day's two quality-grade records).
$ws.Rows("38:39").Insert()

# Row 38: new "Primera" record for 2021-11-19
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44519
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100107
$ws.Range("H38").Value = "Otros"
$ws.Range("I38").Value = 100107002
$ws.Range("J38").Value = "Chirimoya"
$ws.Range("K38").Value = "Cultivar IV Región"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = 2700
$ws.Range("O38").Value = 3000
$ws.Range("P38").Value = 2850
$ws.Range("Q38").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R38").Value = "Provincia del Elquí"
$ws.Range("S38").Value = 2850
$ws.Range("T38").Value = 1

# Row 39: new "Segunda" record for 2021-11-19
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44519
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100107
$ws.Range("H39").Value = "Otros"
$ws.Range("I39").Value = 100107002
$ws.Range("J39").Value = "Chirimoya"
$ws.Range("K39").Value = "Cultivar IV Región"
$ws.Range("L39").Value = "Segunda"
$ws.Range("M39").Value = 60
$ws.Range("N39").Value = 2500
$ws.Range("O39").Value = 2500
$ws.Range("P39").Value = 2500
$ws.Range("Q39").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R39").Value = "Provincia del Elquí"
$ws.Range("S39").Value = 2500
$ws.Range("T39").Value = 1
